# Apply the two changes described by the diff:
#   1. "21 years" -> "15+ years" in the PROFESSIONAL SUMMARY paragraph.
#   2. Remove the entire EDUCATION section (the "EDUCATION" Heading2
#      paragraph plus its two Heading3 degree entries).

$d = $word.ActiveDocument

# --- 1. Update the years of experience -------------------------------
$d.Content.Find.Execute(
    "21 years of experience", $true, $false, $false, $false, $false,
    $true, 1, $false, "15+ years of experience", 2) | Out-Null

# --- 2. Remove the EDUCATION section -----------------------------------
# Locate the start of the "EDUCATION" heading paragraph.
$eduRange = $d.Content
$eduRange.Find.Execute(
    "EDUCATION", $true, $true, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$eduStart = $eduRange.Start

# Locate the end of the last EDUCATION entry ("Bachelor of Arts..."),
# expanding to the full paragraph so the trailing paragraph mark is
# included too (otherwise a multi-paragraph delete that stops mid-way
# through the final paragraph's text is a no-op).
$bachRange = $d.Content
$bachRange.Find.Execute(
    "Bachelor of Arts in Political Science - University of California, Berkeley",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$bachRange.Expand(4) | Out-Null    # wdParagraph = 4
$bachEnd = $bachRange.End

$d.Range($eduStart, $bachEnd).Delete()
